$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.161.92"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.328.92"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'304.82"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").Value = "'97.70"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "'35.61"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "'19.66"
$ws.Range("E11").Value = "  +8.29%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'6.96"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "2.691.34"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "2.328.87"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "'0.789"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "43.050.74"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "'12.64"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'67.98"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'237.68"
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'24.97"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.07"
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'165.65"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'9.15"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'33.29"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'18.13"
$ws.Range("E33").Value = "  +6.35%  "
$ws.Range("D34").Value = "'5.01"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "'4.56"
$ws.Range("E35").Value = "  -8.59%  "
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").Value = "'1.77"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.110"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "1.995.29"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'10.77"
$ws.Range("E43").Value = "  +6.45%  "
$ws.Range("D44").Value = "'0.0281"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "'18.10"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "2.557.81"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'53.75"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'72.08"
$ws.Range("E51").Value = "  -0.51%  "
